$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Secciones" column (B) for rows 10-26 with the new section
# codes. Cell values are written in the order the new labels first appear
# (1ESOA, BC1C, BC1B, 2ESOC) so the shared-string table grows in that order.
$ws.Range("B20").Value = "1ESOA"
$ws.Range("B21").Value = "1ESOA"
$ws.Range("B22").Value = "1ESOA"

$ws.Range("B14").Value = "BC1C"
$ws.Range("B15").Value = "BC1C"
$ws.Range("B16").Value = "BC1C"
$ws.Range("B23").Value = "BC1C"
$ws.Range("B24").Value = "BC1C"
$ws.Range("B25").Value = "BC1C"
$ws.Range("B26").Value = "BC1C"

$ws.Range("B17").Value = "BC1B"
$ws.Range("B18").Value = "BC1B"
$ws.Range("B19").Value = "BC1B"

$ws.Range("B10").Value = "2ESOC"
$ws.Range("B11").Value = "2ESOC"
$ws.Range("B12").Value = "2ESOC"
$ws.Range("B13").Value = "2ESOC"

# Restore the on-screen selection to what the author left it at.
$ws.Range("B10:B13").Select()
